$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.032.54"
$ws.Range("E2").Value = "  +1.89%  "
$ws.Range("D3").Value = "1.908.86"
$ws.Range("E3").Value = "  +2.09%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.007"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.66%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.40"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.32%  "
$ws.Range("E6").Value = "  -0.67%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4831"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.91%  "
$ws.Range("E8").Value = "  +0.91%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07360"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -0.02%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9351"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -0.27%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.82"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.54%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07774"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -1.00%  "
$ws.Range("D13").Value = "1.914.01"
$ws.Range("E13").Value = "  +2.50%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.503"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +1.08%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.628"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.84%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "91.86"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +1.23%  "
$ws.Range("E17").Value = "  -0.73%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008835"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.90%  "
$ws.Range("E19").Value = "  -0.62%  "
$ws.Range("D20").Value = "28.072.90"
$ws.Range("E20").Value = "  +2.04%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.81"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.81%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.182"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +1.06%  "
$ws.Range("D23").Value = "2.162.35"
$ws.Range("E23").Value = "  +3.57%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.92"
$ws.Range("D24").ClearFormats()
$ws.Range("B25").Value = "Monero"
$ws.Range("C25").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "155.75"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +1.26%  "
$ws.Range("B26").Value = "Toncoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.923"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -1.39%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.51"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.04%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.126"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +4.75%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "116.69"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +0.72%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.974"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.37%  "
$ws.Range("E31").Value = "  +0.43%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.304"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -0.89%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.252"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +2.99%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7754"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +2.71%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.679"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +1.51%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.651"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -2.18%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02058"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +0.45%  "
$ws.Range("E38").Value = "  -0.64%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05305"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +0.58%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5485"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +2.58%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.995"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.17%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "7.030"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.78%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.517"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +0.11%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1528"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +0.16%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.72"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.03%  "
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4831"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.44%  "
$ws.Range("B47").Value = "Quant"
$ws.Range("C47").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "108.24"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +5.18%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.006"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.72%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.651"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.66%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "67.95"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +0.78%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06081"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.07%  "
